# Auto-generated edit script: updates cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Sets $value into the cell at $cellRef as a literal text string, even when
    # $value looks like a number (avoids Excel auto-converting it to a numeric cell).
    # Builds the value as a quoted-string formula, then pastes-special as values only,
    # so the result is a plain literal (no residual formula, no style changes).
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $escapedValue = $value -replace '"', '""'
    $r.Formula = '="' + $escapedValue + '"'
    $r.Copy($r)
    $r.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# --- Direct text/string assignments (safe from numeric auto-coercion) ---
$ws.Range("D2").Value = "35.141.03"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.893.34"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("E6").Value = "  +5.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -4.12%  "
$ws.Range("E9").Value = "  +4.19%  "
$ws.Range("E10").Value = "  +13.13%  "
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "2.168.64"
$ws.Range("E14").Value = "  +5.59%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.898.71"
$ws.Range("E16").Value = "  +1.54%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "35.125.96"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").Value = "0.0$([char]0x2083)0815"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +25.14%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("E30").Value = "  +1.23%  "
$ws.Range("D31").Value = "4.162.80"
$ws.Range("E31").Value = "  +21.93%  "
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("E33").Value = "  +14.31%  "
$ws.Range("E34").Value = "  -0.77%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +1.97%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E42").Value = "  +5.35%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "1.331.88"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  +38.88%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "2.078.27"

# --- Numeric-looking text assignments (use Set-TextValue to force literal text) ---
$cellRef = "D5"
$cellValue = "245.54"
Set-TextValue $cellRef $cellValue
$cellRef = "D6"
$cellValue = "0.659"
Set-TextValue $cellRef $cellValue
$cellRef = "D8"
$cellValue = "41.13"
Set-TextValue $cellRef $cellValue
$cellRef = "D9"
$cellValue = "0.345"
Set-TextValue $cellRef $cellValue
$cellRef = "D10"
$cellValue = "52.94"
Set-TextValue $cellRef $cellValue
$cellRef = "D11"
$cellValue = "0.0717"
Set-TextValue $cellRef $cellValue
$cellRef = "D12"
$cellValue = "0.0991"
Set-TextValue $cellRef $cellValue
$cellRef = "D14"
$cellValue = "12.23"
Set-TextValue $cellRef $cellValue
$cellRef = "D15"
$cellValue = "0.694"
Set-TextValue $cellRef $cellValue
$cellRef = "D17"
$cellValue = "4.78"
Set-TextValue $cellRef $cellValue
$cellRef = "D19"
$cellValue = "71.82"
Set-TextValue $cellRef $cellValue
$cellRef = "D21"
$cellValue = "239.96"
Set-TextValue $cellRef $cellValue
$cellRef = "D23"
$cellValue = "4.79"
Set-TextValue $cellRef $cellValue
$cellRef = "D27"
$cellValue = "170.41"
Set-TextValue $cellRef $cellValue
$cellRef = "D28"
$cellValue = "8.41"
Set-TextValue $cellRef $cellValue
$cellRef = "D29"
$cellValue = "18.27"
Set-TextValue $cellRef $cellValue
$cellRef = "D30"
$cellValue = "0.126"
Set-TextValue $cellRef $cellValue
$cellRef = "D33"
$cellValue = "0.939"
Set-TextValue $cellRef $cellValue
$cellRef = "D34"
$cellValue = "0.0559"
Set-TextValue $cellRef $cellValue
$cellRef = "D36"
$cellValue = "4.07"
Set-TextValue $cellRef $cellValue
$cellRef = "D38"
$cellValue = "2.01"
Set-TextValue $cellRef $cellValue
$cellRef = "D40"
$cellValue = "1.09"
Set-TextValue $cellRef $cellValue
$cellRef = "D42"
$cellValue = "16.01"
Set-TextValue $cellRef $cellValue
$cellRef = "D43"
$cellValue = "0.0633"
Set-TextValue $cellRef $cellValue
$cellRef = "D44"
$cellValue = "89.37"
Set-TextValue $cellRef $cellValue
$cellRef = "D46"
$cellValue = "48.65"
Set-TextValue $cellRef $cellValue
$cellRef = "D47"
$cellValue = "2.35"
Set-TextValue $cellRef $cellValue
